$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they reuse the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

# Set the header text/values for the new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells (plain/default style, matching H2)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9
